$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(29).Insert()
$ws.Range("A30:Q30").Copy()
$ws.Range("A29:Q29").PasteSpecial(-4122)
$ws.Rows.Item(29).RowHeight = 25.5
